$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates from the cryptos list refresh.
# Numeric-looking Price values must be forced to Text format so they
# keep matching the original inline-string/text representation instead
# of being auto-converted to numbers by Excel.

$ws.Range('D2').Value = '64.377.72'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '3.135.19'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '572.18'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.53'
$ws.Range('E6').Value = '  -3.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -6.49%  '
$ws.Range('D9').Value = '3.143.81'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.61'
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.382'
$ws.Range('E12').Value = '  -3.06%  '
$ws.Range('D13').Value = '3.675.76'
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '64.435.05'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '24.81'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').Value = '3.134.84'
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '411.18'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.21'
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.46'
$ws.Range('E21').Value = '  -4.26%  '
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.51'
$ws.Range('E24').Value = '  -2.61%  '
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.192'
$ws.Range('E26').Value = '  -6.74%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000102'
$ws.Range('E27').Value = '  -3.84%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.85'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -2.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '21.20'
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '162.58'
$ws.Range('E33').Value = '  +4.24%  '
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.23'
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  -2.54%  '
$ws.Range('D39').Value = '2.626.60'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.10'
$ws.Range('E40').Value = '  -3.93%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '23.53'
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '38.17'
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.687'
$ws.Range('E43').Value = '  -4.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0612'
$ws.Range('E44').Value = '  -2.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.29'
$ws.Range('E45').Value = '  -4.80%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '289.04'
$ws.Range('E46').Value = '  -1.66%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '21.25'
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0254'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.995'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0973'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.48'
$ws.Range('E51').Value = '  +0.45%  '
